$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5527
$ws1.Range("F9").Value = 524

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5527
$ws4.Range("F11").Value = 524
